$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the existing last row (row 6: 7 / 9) down to row 8, then insert the
# two new rows (5/0 and 6/0) at rows 6 and 7.
$lastA = $ws.Range("A6").Value2
$lastB = $ws.Range("B6").Value2
$ws.Range("A8").Value2 = $lastA
$ws.Range("B8").Value2 = $lastB

$ws.Range("A6").Value2 = 5
$ws.Range("B6").Value2 = 0

$ws.Range("A7").Value2 = 6
$ws.Range("B7").Value2 = 0

# Update the active cell / selection to match the post-import cursor position.
$ws.Range("C9").Select()
